$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111271309
$ws.Range("B2").Value = 78579
$ws.Range("E2").Value = 2081
$ws.Range("F2").Value = "Skrovellav"
$ws.Range("G2").Value = "Lobaria scrobiculata"
$ws.Range("H2").Value = "(Scop.) DC."
$ws.Range("Q2").Value = 473221.4734201821
$ws.Range("R2").Value = 6863586.84377678

# Row 3
$ws.Range("A3").Value = 111270596
$ws.Range("Q3").Value = 473184.8241620373
$ws.Range("R3").Value = 6863788.37406126
$ws.Range("AC3").Value = "Fem blommande."

# Row 4
$ws.Range("A4").Value = 111271176
$ws.Range("B4").Value = 78579
$ws.Range("E4").Value = 2081
$ws.Range("F4").Value = "Skrovellav"
$ws.Range("G4").Value = "Lobaria scrobiculata"
$ws.Range("H4").Value = "(Scop.) DC."
$ws.Range("Q4").Value = 473227.9160841404
$ws.Range("R4").Value = 6863625.911539786

# Row 5
$ws.Range("A5").Value = 111270755
$ws.Range("Q5").Value = 473194.7999623233
$ws.Range("R5").Value = 6863736.454484907

# Row 6
$ws.Range("A6").Value = 111271821
$ws.Range("AC6").Value = "Tre blommande."

# Row 7
$ws.Range("A7").Value = 111271068
$ws.Range("B7").Value = 78578
$ws.Range("E7").Value = 6458
$ws.Range("F7").Value = "Lunglav"
$ws.Range("G7").Value = "Lobaria pulmonaria"
$ws.Range("H7").Value = "(L.) Hoffm."
$ws.Range("Q7").Value = 473238.8676645419
$ws.Range("R7").Value = 6863638.079474191

# Row 8
$ws.Range("A8").Value = 111271061
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."

# Row 9
$ws.Range("A9").Value = 111272292
$ws.Range("B9").Value = 78578
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6458
$ws.Range("F9").Value = "Lunglav"
$ws.Range("G9").Value = "Lobaria pulmonaria"
$ws.Range("H9").Value = "(L.) Hoffm."
$ws.Range("Q9").Value = 473321.1690919191
$ws.Range("R9").Value = 6863539.403128584

# Row 10
$ws.Range("A10").Value = 111271382
$ws.Range("Q10").Value = 473167.6377000402
$ws.Range("R10").Value = 6863583.496200636
$ws.Range("AC10").Value = "Tre blommande."

# Row 11
$ws.Range("A11").Value = 111271055
$ws.Range("Q11").Value = 473238.8676645419
$ws.Range("R11").Value = 6863638.079474191

# Row 12
$ws.Range("A12").Value = 111270939
$ws.Range("B12").Value = 78578
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6458
$ws.Range("F12").Value = "Lunglav"
$ws.Range("G12").Value = "Lobaria pulmonaria"
$ws.Range("H12").Value = "(L.) Hoffm."
$ws.Range("Q12").Value = 473229.5908188519
$ws.Range("R12").Value = 6863658.889402787
$ws.Range("AC12").ClearContents() | Out-Null

# Row 15
$ws.Range("A15").Value = 111270892
$ws.Range("B15").Value = 78578
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 6458
$ws.Range("F15").Value = "Lunglav"
$ws.Range("G15").Value = "Lobaria pulmonaria"
$ws.Range("H15").Value = "(L.) Hoffm."
$ws.Range("Q15").Value = 473239.9383552746
$ws.Range("R15").Value = 6863714.420922431
$ws.Range("AC15").ClearContents() | Out-Null

# Row 16
$ws.Range("A16").Value = 111271296
$ws.Range("Q16").Value = 473220.1559155915
$ws.Range("R16").Value = 6863539.25170773

# Row 18
$ws.Range("A18").Value = 111270747
$ws.Range("B18").Value = 96348
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 220787
$ws.Range("F18").Value = "Knärot"
$ws.Range("G18").Value = "Goodyera repens"
$ws.Range("H18").Value = "(L.) R. Br."
$ws.Range("AC18").Value = "Sex blommande."

# Row 19
$ws.Range("A19").Value = 111271923
$ws.Range("B19").Value = 96348
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("Q19").Value = 473118.5439814709
$ws.Range("R19").Value = 6863582.939962601
$ws.Range("AC19").Value = "Tre blommande."

# Row 20
$ws.Range("A20").Value = 111271029
$ws.Range("B20").Value = 78579
$ws.Range("E20").Value = 2081
$ws.Range("F20").Value = "Skrovellav"
$ws.Range("G20").Value = "Lobaria scrobiculata"
$ws.Range("H20").Value = "(Scop.) DC."
$ws.Range("Q20").Value = 473229.5908188519
$ws.Range("R20").Value = 6863658.889402787

# Row 21
$ws.Range("A21").Value = 111270559
$ws.Range("B21").Value = 96348
$ws.Range("D21").Value = "VU"
$ws.Range("E21").Value = 220787
$ws.Range("F21").Value = "Knärot"
$ws.Range("G21").Value = "Goodyera repens"
$ws.Range("H21").Value = "(L.) R. Br."
$ws.Range("Q21").Value = 473167.8634183492
$ws.Range("R21").Value = 6863792.277629613
$ws.Range("AC21").Value = "Två blommande."

# Row 22
$ws.Range("A22").Value = 111271588
$ws.Range("B22").Value = 78578
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6458
$ws.Range("F22").Value = "Lunglav"
$ws.Range("G22").Value = "Lobaria pulmonaria"
$ws.Range("H22").Value = "(L.) Hoffm."
$ws.Range("Q22").Value = 473140.3516782348
$ws.Range("R22").Value = 6863595.022241795
$ws.Range("AC22").ClearContents() | Out-Null
